$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$updates = @(
    @{ Row = 3;   C = 249326;  D = $null; E = 1036473141 }
    @{ Row = 8;   C = 1045;    D = 194;   E = 91109976 }
    @{ Row = 36;  C = 56967;   D = $null; E = 223721899 }
    @{ Row = 53;  C = 141678;  D = $null; E = 590056680 }
    @{ Row = 54;  C = 60292;   D = $null; E = 353611523 }
    @{ Row = 62;  C = 4188;    D = $null; E = 9183365 }
    @{ Row = 63;  C = 14352;   D = $null; E = 36179469 }
    @{ Row = 64;  C = 5198;    D = $null; E = 20337729 }
    @{ Row = 70;  C = 15724;   D = $null; E = 24658131 }
    @{ Row = 91;  C = 151093;  D = $null; E = 482071222 }
    @{ Row = 92;  C = 408994;  D = $null; E = 1593538885 }
    @{ Row = 93;  C = 209480;  D = $null; E = 1307421313 }
    @{ Row = 94;  C = 94140;   D = $null; E = 915436245 }
    @{ Row = 95;  C = 50709;   D = $null; E = 929737465 }
    @{ Row = 96;  C = 17240;   D = $null; E = 789174173 }
    @{ Row = 98;  C = 809;     D = $null; E = 117674774 }
    @{ Row = 107; C = 6387;    D = $null; E = 21941144 }
    @{ Row = 109; C = 1271;    D = $null; E = 20757709 }
    @{ Row = 114; C = 3791;    D = $null; E = 9080369 }
    @{ Row = 115; C = 11691;   D = $null; E = 32947448 }
    @{ Row = 116; C = 4554;    D = $null; E = 20454780 }
    @{ Row = 118; C = 973;     D = $null; E = 11742670 }
    @{ Row = 122; C = 8485;    D = $null; E = 12669783 }
    @{ Row = 124; C = 948;     D = $null; E = 2624236 }
    @{ Row = 142; C = 168972;  D = $null; E = 681760934 }
    @{ Row = 156; C = 25104;   D = $null; E = 199245925 }
)

foreach ($u in $updates) {
    $r = $u.Row
    $ws.Range("C$r").Value = $u.C
    if ($null -ne $u.D) {
        $ws.Range("D$r").Value = $u.D
    }
    $ws.Range("E$r").Value = $u.E
}
